$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.213846683502197
$ws.Range("B1").Value = 2.555657625198364
$ws.Range("C1").Value = 9.103615760803223
$ws.Range("D1").Value = 2.071223735809326
$ws.Range("E1").Value = 1.19152569770813
